$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Escopo")
$ws.Activate()

# Update F4 (row 4) percentage from 0.2 to 0.3
$ws.Range("F4").Value = 0.3

# E9: set Status to "Em Andamento" and F9 from 0 to 0.15
$ws.Range("E9").Value = "Em Andamento"
$ws.Range("F9").Value = 0.15

# E10: set Status to "Em Andamento" and F10 from 0 to 0.15
$ws.Range("E10").Value = "Em Andamento"
$ws.Range("F10").Value = 0.15

# E12: set Status to "Em Andamento" (F12 value stays 0.25)
$ws.Range("E12").Value = "Em Andamento"

# E23: set Status to "Em Andamento" and F23 from 0 to 0.5
$ws.Range("E23").Value = "Em Andamento"
$ws.Range("F23").Value = 0.5

# Select B1 as active cell on this sheet (matches sheetView selection change)
$ws.Range("B1").Select()
